# Auto-generated Excel COM-interop script applying numeric corrections
# to the LeveProfits tables across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 12197468
$ws.Range("J112").Value = 12822714
$ws.Range("L112").Value = 38468142
$ws.Range("N112").Value = -38470358

$ws.Range("H113").Value = 7067.52
$ws.Range("J113").Value = 11709.444
$ws.Range("L113").Value = 11709.444
$ws.Range("N113").Value = -18217.444

$ws.Range("H127").Value = 1447.5625
$ws.Range("I127").Value = 709.5
$ws.Range("J127").Value = 3661.75
$ws.Range("K127").Value = 2128.5
$ws.Range("L127").Value = 10985.25
$ws.Range("M127").Value = 2831.5
$ws.Range("N127").Value = -20905.25

$ws.Range("H132").Value = 18692.143
$ws.Range("I132").Value = 18983.846
$ws.Range("K132").Value = 56951.538
$ws.Range("M132").Value = -54421.538

$ws.Range("H137").Value = 7852.222
$ws.Range("I137").Value = 4942.4707
$ws.Range("J137").Value = 12798.8
$ws.Range("K137").Value = 14827.4121
$ws.Range("L137").Value = 38396.39999999999
$ws.Range("M137").Value = -12277.4121
$ws.Range("N137").Value = -43496.39999999999

$ws.Range("H138").Value = 20003160
$ws.Range("J138").Value = 33337656
$ws.Range("L138").Value = 100012968
$ws.Range("N138").Value = -100023248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 585.5700000000001
$ws.Range("I32").Value = 579.1739
$ws.Range("J32").Value = 659.125
$ws.Range("K32").Value = 579.1739
$ws.Range("L32").Value = 659.125
$ws.Range("M32").Value = -292.1739
$ws.Range("N32").Value = -1233.125

$ws.Range("H61").Value = 3882.6365
$ws.Range("I61").Value = 2260.1
$ws.Range("J61").Value = 5234.75
$ws.Range("K61").Value = 2260.1
$ws.Range("L61").Value = 5234.75
$ws.Range("M61").Value = -2048.1
$ws.Range("N61").Value = -5658.75

$ws.Range("H110").Value = 33388.184
$ws.Range("I110").Value = 44761.875
$ws.Range("J110").Value = 3058.3333
$ws.Range("K110").Value = 44761.875
$ws.Range("L110").Value = 3058.3333
$ws.Range("M110").Value = -42716.875
$ws.Range("N110").Value = -7148.3333

$ws.Range("H136").Value = 3882.6365
$ws.Range("I136").Value = 2260.1
$ws.Range("J136").Value = 5234.75
$ws.Range("K136").Value = 6780.299999999999
$ws.Range("L136").Value = 15704.25
$ws.Range("M136").Value = -4230.299999999999
$ws.Range("N136").Value = -20804.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1897.9048
$ws.Range("I20").Value = 2788.8
$ws.Range("J20").Value = 1088
$ws.Range("K20").Value = 2788.8
$ws.Range("L20").Value = 1088
$ws.Range("M20").Value = -2541.8
$ws.Range("N20").Value = -1582

$ws.Range("H26").Value = 15623.363
$ws.Range("I26").Value = 15623.363
$ws.Range("K26").Value = 15623.363
$ws.Range("M26").Value = -15331.363

$ws.Range("H107").Value = 2924.6
$ws.Range("I107").Value = 2351.9333
$ws.Range("K107").Value = 2351.9333
$ws.Range("M107").Value = -431.9333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1187.625
$ws.Range("I16").Value = 1231.5
$ws.Range("K16").Value = 1231.5
$ws.Range("M16").Value = -944.5

$ws.Range("H31").Value = 6681.6924
$ws.Range("I31").Value = 4909
$ws.Range("J31").Value = 7981.6665
$ws.Range("K31").Value = 4909
$ws.Range("L31").Value = 7981.6665
$ws.Range("M31").Value = -4614
$ws.Range("N31").Value = -8571.666499999999

$ws.Range("H34").Value = 6681.6924
$ws.Range("I34").Value = 4909
$ws.Range("J34").Value = 7981.6665
$ws.Range("K34").Value = 4909
$ws.Range("L34").Value = 7981.6665
$ws.Range("M34").Value = -4707
$ws.Range("N34").Value = -8385.666499999999

$ws.Range("H57").Value = 39996.5

$ws.Range("H58").Value = 3860.8147
$ws.Range("I58").Value = 3835.875
$ws.Range("J58").Value = 3897.0908
$ws.Range("K58").Value = 3835.875
$ws.Range("L58").Value = 3897.0908
$ws.Range("M58").Value = -3632.875
$ws.Range("N58").Value = -4303.0908

$ws.Range("H107").Value = 527
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 527
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 527
$ws.Range("N107").Value = -4367
$ws.Range("M107").ClearContents()

$ws.Range("H113").Value = 1187.625
$ws.Range("I113").Value = 1231.5
$ws.Range("K113").Value = 1231.5
$ws.Range("M113").Value = 938.5

$ws.Range("H132").Value = 116917.664
$ws.Range("I132").Value = 148402.08
$ws.Range("K132").Value = 445206.24
$ws.Range("M132").Value = -442676.24

$ws.Range("H134").Value = 4029.5
$ws.Range("I134").Value = 3374.077
$ws.Range("K134").Value = 10122.231
$ws.Range("M134").Value = -7587.231

$ws.Range("H136").Value = 3860.8147
$ws.Range("I136").Value = 3835.875
$ws.Range("J136").Value = 3897.0908
$ws.Range("K136").Value = 11507.625
$ws.Range("L136").Value = 11691.2724
$ws.Range("M136").Value = -8957.625
$ws.Range("N136").Value = -16791.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1765.2222
$ws.Range("J47").Value = 2291.6667
$ws.Range("L47").Value = 6875.000100000001
$ws.Range("N47").Value = -7737.000100000001

$ws.Range("H123").Value = 1660.8

$ws.Range("H125").Value = 4016.25
$ws.Range("I125").Value = 3999.5
$ws.Range("J125").Value = 4033
$ws.Range("K125").Value = 11998.5
$ws.Range("L125").Value = 12099
$ws.Range("M125").Value = -7078.5
$ws.Range("N125").Value = -21939

$ws.Range("H131").Value = 1420.4642
$ws.Range("J131").Value = 1541.8948
$ws.Range("L131").Value = 4625.6844
$ws.Range("N131").Value = -14705.6844

$ws.Range("H132").Value = 3544.2856
$ws.Range("I132").Value = 1665.2858
$ws.Range("J132").Value = 5423.2856
$ws.Range("K132").Value = 14987.5722
$ws.Range("L132").Value = 48809.5704
$ws.Range("M132").Value = -12457.5722
$ws.Range("N132").Value = -53869.5704

$ws.Range("H134").Value = 3751.5
$ws.Range("I134").Value = 1155.7693
$ws.Range("J134").Value = 14999.667
$ws.Range("K134").Value = 3467.3079
$ws.Range("L134").Value = 44999.001
$ws.Range("M134").Value = 1602.6921
$ws.Range("N134").Value = -55139.001

$ws.Range("H137").Value = 1099.2
$ws.Range("I137").Value = 1099.2
$ws.Range("K137").Value = 3297.6
$ws.Range("M137").Value = 1802.4

$ws.Range("H139").Value = 2174.6428
$ws.Range("I139").Value = 2031.5834
$ws.Range("K139").Value = 6094.7502
$ws.Range("M139").Value = -954.7502000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 202.625
$ws.Range("I2").Value = 51.466667
$ws.Range("J2").Value = 336
$ws.Range("K2").Value = 51.466667
$ws.Range("L2").Value = 336
$ws.Range("M2").Value = 61.533333
$ws.Range("N2").Value = -562

$ws.Range("H49").Value = 26958.166
$ws.Range("J49").Value = 24499.5
$ws.Range("L49").Value = 24499.5
$ws.Range("N49").Value = -24867.5

$ws.Range("H126").Value = 2124.6904
$ws.Range("I126").Value = 2030.0625
$ws.Range("J126").Value = 2427.5
$ws.Range("K126").Value = 6090.1875
$ws.Range("L126").Value = 7282.5
$ws.Range("M126").Value = -3620.1875
$ws.Range("N126").Value = -12222.5

$ws.Range("H132").Value = 1338.3334
$ws.Range("I132").Value = 1278.1818
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3834.5454
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1304.5454
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1453.7916
$ws.Range("I22").Value = 1046.5714
$ws.Range("J22").Value = 2023.9
$ws.Range("K22").Value = 1046.5714
$ws.Range("L22").Value = 2023.9
$ws.Range("M22").Value = -751.5714
$ws.Range("N22").Value = -2613.9

$ws.Range("H27").Value = 1453.7916
$ws.Range("I27").Value = 1046.5714
$ws.Range("J27").Value = 2023.9
$ws.Range("K27").Value = 1046.5714
$ws.Range("L27").Value = 2023.9
$ws.Range("M27").Value = -939.5714
$ws.Range("N27").Value = -2237.9

$ws.Range("H46").Value = 1414.2858

$ws.Range("H48").Value = 25589.834
$ws.Range("I48").Value = 23497.5
$ws.Range("J48").Value = 29774.5
$ws.Range("K48").Value = 23497.5
$ws.Range("L48").Value = 29774.5
$ws.Range("M48").Value = -22836.5
$ws.Range("N48").Value = -31096.5

$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180

$ws.Range("H122").Value = 3898.1667
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 4694.5
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 14083.5
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -18983.5

$ws.Range("H132").Value = 8733.735000000001
$ws.Range("I132").Value = 2998.1072
$ws.Range("K132").Value = 8994.321599999999
$ws.Range("M132").Value = -6464.321599999999

$ws.Range("H136").Value = 2263.1943
$ws.Range("I136").Value = 1824.9032
$ws.Range("K136").Value = 5474.7096
$ws.Range("M136").Value = -2924.7096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 33332.668
$ws.Range("J47").Value = 33999
$ws.Range("L47").Value = 33999
$ws.Range("N47").Value = -35143

$ws.Range("H52").Value = 22421.834
$ws.Range("J52").Value = 27745
$ws.Range("L52").Value = 27745
$ws.Range("N52").Value = -28197

$ws.Range("H70").Value = 79999
$ws.Range("J70").Value = 79999
$ws.Range("L70").Value = 79999
$ws.Range("N70").Value = -80629

$ws.Range("H73").Value = 79999
$ws.Range("J73").Value = 79999
$ws.Range("L73").Value = 79999
$ws.Range("N73").Value = -82183

$ws.Range("H113").Value = 1379.0625
$ws.Range("I113").Value = 823.63635
$ws.Range("J113").Value = 2601
$ws.Range("K113").Value = 2470.90905
$ws.Range("L113").Value = 7803
$ws.Range("M113").Value = -300.9090500000002
$ws.Range("N113").Value = -12143

$ws.Range("H136").Value = 5713.231
$ws.Range("I136").Value = 4474.085
$ws.Range("J136").Value = 9561.105
$ws.Range("K136").Value = 13422.255
$ws.Range("L136").Value = 28683.315
$ws.Range("M136").Value = -10872.255
$ws.Range("N136").Value = -33783.315
